$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-26 Sunday" "2024-05-27 Monday"

Replace-Text "827×3=" "345×4="
Replace-Text "581×6=" "154×3="
Replace-Text "247×4=" "859×7="
Replace-Text "749×8=" "223×8="
Replace-Text "717×5=" "760×2="

Replace-Text "768×8=" "163×6="
Replace-Text "633×4=" "226×9="
Replace-Text "371×9=" "535×8="
Replace-Text "918×8=" "676×7="
Replace-Text "591×9=" "174×2="

Replace-Text "466×6=" "453×8="
Replace-Text "330×2=" "172×9="
Replace-Text "907×4=" "300×8="
Replace-Text "935×6=" "280×6="
Replace-Text "456×2=" "536×7="

Replace-Text "878×5=" "236×3="
Replace-Text "581×2=" "505×9="
Replace-Text "115×7=" "728×8="
Replace-Text "557×8=" "481×8="
Replace-Text "650×3=" "564×4="

Replace-Text "689×8=" "735×9="
Replace-Text "306×9=" "410×9="
Replace-Text "876×5=" "578×4="
Replace-Text "416×3=" "196×4="
Replace-Text "301×8=" "494×8="
